$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 8
$ws.Range("H8").Value = 41
$ws.Range("I8").Value = 30.272728
$ws.Range("J8").Value = 100
$ws.Range("K8").Value = 90.818184
$ws.Range("L8").Value = 300
$ws.Range("M8").Value = 48.181816
$ws.Range("N8").Value = -578
# Row 62
$ws.Range("H62").Value = 6712.857
$ws.Range("J62").Value = 8497.5
$ws.Range("L62").Value = 8497.5
$ws.Range("N62").Value = -9745.5
# Row 65
$ws.Range("H65").Value = 6712.857
$ws.Range("J65").Value = 8497.5
$ws.Range("L65").Value = 42487.5
$ws.Range("N65").Value = -48727.5
# Row 106
$ws.Range("H106").Value = 8600
$ws.Range("I106").Value = 8600
$ws.Range("K106").Value = 8600
$ws.Range("M106").Value = -7969
# Row 112
$ws.Range("H112").Value = 2198.6667
$ws.Range("J112").Value = 2398.1667
$ws.Range("L112").Value = 7194.500100000001
$ws.Range("N112").Value = -9410.500100000001
# Row 137
$ws.Range("H137").Value = 1948.4667
$ws.Range("I137").Value = 1792.7
$ws.Range("J137").Value = 2260
$ws.Range("K137").Value = 5378.1
$ws.Range("L137").Value = 6780
$ws.Range("M137").Value = -2828.1
$ws.Range("N137").Value = -11880

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 2534.1333
$ws.Range("I2").Value = 1833.2222
$ws.Range("K2").Value = 1833.2222
$ws.Range("M2").Value = -1720.2222
# Row 23
$ws.Range("H23").Value = 20000
$ws.Range("I23").Value = 20000
$ws.Range("K23").Value = 20000
$ws.Range("M23").Value = -19741
# Row 37
$ws.Range("H37").Value = 23641.428
$ws.Range("J37").Value = 23641.428
$ws.Range("L37").Value = 23641.428
$ws.Range("N37").Value = -24187.428
# Row 44
$ws.Range("H44").Value = 34996.668
$ws.Range("J44").Value = 34996.668
$ws.Range("L44").Value = 34996.668
$ws.Range("N44").Value = -35972.668
# Row 46
$ws.Range("H46").Value = 8427.2
$ws.Range("I46").Value = 7712.6665
$ws.Range("J46").Value = 9499
$ws.Range("K46").Value = 7712.6665
$ws.Range("L46").Value = 9499
$ws.Range("M46").Value = -7393.6665
$ws.Range("N46").Value = -10137
# Row 55
$ws.Range("H55").Value = 24998.75
$ws.Range("J55").Value = 24998.75
$ws.Range("L55").Value = 24998.75
$ws.Range("N55").Value = -25628.75
# Row 61
$ws.Range("H61").Value = 4200.143
$ws.Range("I61").Value = 2999.6667
$ws.Range("J61").Value = 5100.5
$ws.Range("K61").Value = 2999.6667
$ws.Range("L61").Value = 5100.5
$ws.Range("M61").Value = -2787.6667
$ws.Range("N61").Value = -5524.5
# Row 74
$ws.Range("H74").Value = 1930.6666
$ws.Range("I74").Value = 1930.6666
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 1930.6666
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -1056.6666
$ws.Range("N74").Value = $null
# Row 77
$ws.Range("H77").Value = 1930.6666
$ws.Range("I77").Value = 1930.6666
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 9653.333
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -5285.333000000001
$ws.Range("N77").Value = $null
# Row 80
$ws.Range("H80").Value = 40000
$ws.Range("J80").Value = 40000
$ws.Range("L80").Value = 40000
$ws.Range("N80").Value = -41996
# Row 83
$ws.Range("H83").Value = 40000
$ws.Range("J83").Value = 40000
$ws.Range("L83").Value = 120000
$ws.Range("N83").Value = -129984
# Row 116
$ws.Range("H116").Value = 2534.1333
$ws.Range("I116").Value = 1833.2222
$ws.Range("K116").Value = 1833.2222
$ws.Range("M116").Value = 460.7778000000001
# Row 136
$ws.Range("H136").Value = 4200.143
$ws.Range("I136").Value = 2999.6667
$ws.Range("J136").Value = 5100.5
$ws.Range("K136").Value = 8999.000100000001
$ws.Range("L136").Value = 15301.5
$ws.Range("M136").Value = -6449.000100000001
$ws.Range("N136").Value = -20401.5

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 2534.1333
$ws.Range("I3").Value = 1833.2222
$ws.Range("K3").Value = 1833.2222
$ws.Range("M3").Value = -1719.2222
# Row 35
$ws.Range("H35").Value = 14997.223
$ws.Range("J35").Value = 14997.223
$ws.Range("L35").Value = 14997.223
$ws.Range("N35").Value = -15617.223
# Row 75
$ws.Range("H75").Value = 21648.625
$ws.Range("I75").Value = 7598.4287
$ws.Range("K75").Value = 7598.4287
$ws.Range("M75").Value = -6662.4287
# Row 78
$ws.Range("H78").Value = 21648.625
$ws.Range("I78").Value = 7598.4287
$ws.Range("K78").Value = 22795.2861
$ws.Range("M78").Value = -18115.2861
# Row 82
$ws.Range("H82").Value = 23179.75
$ws.Range("I82").Value = 6360.3335
$ws.Range("J82").Value = 39999.168
$ws.Range("K82").Value = 6360.3335
$ws.Range("L82").Value = 39999.168
$ws.Range("M82").Value = -5977.3335
$ws.Range("N82").Value = -40765.168
# Row 85
$ws.Range("H85").Value = 23179.75
$ws.Range("I85").Value = 6360.3335
$ws.Range("J85").Value = 39999.168
$ws.Range("K85").Value = 6360.3335
$ws.Range("L85").Value = 39999.168
$ws.Range("M85").Value = -5034.3335
$ws.Range("N85").Value = -42651.168
# Row 94
$ws.Range("H94").Value = 1334
$ws.Range("I94").Value = 1141.0526
$ws.Range("K94").Value = 1141.0526
$ws.Range("M94").Value = -690.0526

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 41
$ws.Range("H41").Value = 17490
$ws.Range("I41").Value = 15000
$ws.Range("J41").Value = 17716.363
$ws.Range("K41").Value = 15000
$ws.Range("L41").Value = 17716.363
$ws.Range("M41").Value = -14572
$ws.Range("N41").Value = -18572.363
# Row 50
$ws.Range("H50").Value = 28165.834
$ws.Range("J50").Value = 28165.834
$ws.Range("L50").Value = 28165.834
$ws.Range("N50").Value = -29415.834
# Row 59
$ws.Range("H59").Value = 29260.572
$ws.Range("I59").Value = 23522.572
$ws.Range("K59").Value = 23522.572
$ws.Range("M59").Value = -22377.572
# Row 68
$ws.Range("H68").Value = 39998.57
$ws.Range("J68").Value = 39998.57
$ws.Range("L68").Value = 39998.57
$ws.Range("N68").Value = -41496.57
# Row 71
$ws.Range("H71").Value = 39998.57
$ws.Range("J71").Value = 39998.57
$ws.Range("L71").Value = 119995.71
$ws.Range("N71").Value = -127483.71
# Row 74
$ws.Range("H74").Value = 39997.5
$ws.Range("J74").Value = 39997.5
$ws.Range("L74").Value = 39997.5
$ws.Range("N74").Value = -41745.5
# Row 77
$ws.Range("H77").Value = 39997.5
$ws.Range("J77").Value = 39997.5
$ws.Range("L77").Value = 119992.5
$ws.Range("N77").Value = -128728.5
# Row 92
$ws.Range("H92").Value = 55241.668
$ws.Range("J92").Value = 55241.668
$ws.Range("L92").Value = 55241.668
$ws.Range("N92").Value = -60233.668
# Row 94
$ws.Range("H94").Value = 1699.375
$ws.Range("I94").Value = 1752
$ws.Range("J94").Value = 1646.75
$ws.Range("K94").Value = 1752
$ws.Range("L94").Value = 1646.75
$ws.Range("M94").Value = -1301
$ws.Range("N94").Value = -2548.75

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 14
$ws.Range("H14").Value = 24733.2
$ws.Range("I14").Value = 24733.2
$ws.Range("K14").Value = 74199.6
$ws.Range("M14").Value = -74026.6
# Row 38
$ws.Range("H38").Value = 33.363636
$ws.Range("I38").Value = 32
$ws.Range("K38").Value = 96
$ws.Range("M38").Value = 251
# Row 138
$ws.Range("H138").Value = 638.3333
$ws.Range("I138").Value = 638.3333
$ws.Range("K138").Value = 1914.9999
$ws.Range("M138").Value = 3225.0001

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 696
$ws.Range("I7").Value = 696
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 696
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -584
$ws.Range("N7").Value = $null
# Row 13
$ws.Range("H13").Value = 1000
$ws.Range("I13").Value = 1000
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 1000
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = -860
$ws.Range("N13").Value = $null
# Row 122
$ws.Range("H122").Value = 2979.111
$ws.Range("I122").Value = 2936.5
$ws.Range("K122").Value = 8809.5
$ws.Range("M122").Value = -6359.5
# Row 126
$ws.Range("H126").Value = 696
$ws.Range("I126").Value = 696
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 2088
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = 382
$ws.Range("N126").Value = $null
# Row 132
$ws.Range("H132").Value = 13026.277
$ws.Range("J132").Value = 12333.5
$ws.Range("L132").Value = 37000.5
$ws.Range("N132").Value = -42060.5

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Range("H113").Value = 1681.4445
$ws.Range("I113").Value = 861.8571
$ws.Range("J113").Value = 4550
$ws.Range("K113").Value = 2585.5713
$ws.Range("L113").Value = 13650
$ws.Range("M113").Value = -415.5712999999996
$ws.Range("N113").Value = -17990
# Row 122
$ws.Range("H122").Value = 2387.5
$ws.Range("I122").Value = 2387.5
$ws.Range("K122").Value = 7162.5
$ws.Range("M122").Value = -4712.5
# Row 132
$ws.Range("H132").Value = 3633.1428
$ws.Range("I132").Value = 3086.4
$ws.Range("K132").Value = 9259.2
$ws.Range("M132").Value = -6729.200000000001
